$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.198.56"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "1.861.28"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "

$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4655"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.69%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2822"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06520"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07849"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.15%  "

$ws.Range("D13").Value = "1.868.17"
$ws.Range("E13").Value = "  -1.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.085"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6707"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.65%  "

$ws.Range("D17").Value = "30.194.14"
$ws.Range("E17").Value = "  -0.40%  "

$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.501"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.64%  "

$ws.Range("D21").Value = "2.115.93"
$ws.Range("E21").Value = "  -1.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007264"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.134"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.174"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.915"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.60%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.375"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09688"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.409"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.471"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.067"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04677"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.111"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7023"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.54%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.727"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01845"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.19%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.533"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.236"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.938"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8441"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.40%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.20%  "

$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4147"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.154"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.12%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "935.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.69%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.095"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.72%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05621"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.19%  "
